$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.016.99"
$ws.Range("E2").Value = "  +3.61%  "

$ws.Range("D3").Value = "2.243.54"
$ws.Range("E3").Value = "  +2.13%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.43%  "

$ws.Range("E7").Value = "  +2.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.18%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.71%  "

$ws.Range("E13").Value = "  +2.90%  "

$ws.Range("D14").Value = "2.574.08"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.49%  "

$ws.Range("D16").Value = "2.233.57"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.62%  "

$ws.Range("D18").Value = "43.897.76"
$ws.Range("E18").Value = "  +3.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.03%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.09%  "

$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.22%  "

$ws.Range("E32").Value = "  +10.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0368"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.46%  "

$ws.Range("E38").Value = "  +3.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.21%  "

$ws.Range("E40").Value = "  +23.86%  "

$ws.Range("E41").Value = "  +3.16%  "

$ws.Range("E42").Value = "  +5.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.52%  "

$ws.Range("E44").Value = "  +2.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "

$ws.Range("E46").Value = "  +2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0986"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.448"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +27.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "

# Row 35/36 content swap (Stellar <-> Kaspa)
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.16%  "

$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.123"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.18%  "
